# Add 2022-Q4 data
# -----------------------------------------------------------------------
# 1. Create a new worksheet "2022-Q4" positioned right after "总计" by
#    duplicating the "2022-Q3" sheet (this preserves header styling /
#    borders / bold font exactly), then overwrite its data with the
#    Q4 numbers and trim the extra (now unused) rows.
# 2. Insert a new row into the "总计" (summary) sheet for the 2022-Q4
#    entry, shifting the existing Q3/Q2/Q1 rows down by one, and fix up
#    the running index in column A.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Step 1: build the "2022-Q4" sheet -------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3, $null)          # new copy is placed immediately before 2022-Q3
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Only one fund is held in 2022-Q4, so drop the leftover copied rows 3:6
$q4.Range("A3:H6").Delete()

# Fund-code / name / ratios are stored as TEXT (matches the rest of the
# workbook), so force text format before assigning to avoid Excel's
# automatic "numeric-looking string -> number" coercion (e.g. "009999").
$q4.Range("B2:G2").NumberFormat = "@"

$q4.Cells.Item(2, 1).Value2 = 0
$q4.Cells.Item(2, 2).Value2 = "009999"
$q4.Cells.Item(2, 3).Value2 = "东方中国红利混合"
$q4.Cells.Item(2, 4).Value2 = "0.48"
$q4.Cells.Item(2, 5).Value2 = "90.83"
$q4.Cells.Item(2, 6).Value2 = "5.31"
$q4.Cells.Item(2, 7).Value2 = "0.0255"
$q4.Cells.Item(2, 8).Value2 = 2

# --- Step 2: update the "总计" (summary) sheet ------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows(2).Insert()

# Copy number/cell formatting from the (shifted) old row so the newly
# inserted row matches the plain data-row style instead of Excel's
# default "inherit from row above" behaviour.
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)   # xlPasteFormats

$total.Cells.Item(2, 1).Value2 = 0
$total.Cells.Item(2, 2).Value2 = "2022-Q4"
$total.Cells.Item(2, 3).Value2 = 1
$total.Cells.Item(2, 4).Value2 = 0.03

# Keep the running 0-based index in column A sequential after the insert
$total.Cells.Item(3, 1).Value2 = 1
$total.Cells.Item(4, 1).Value2 = 2
$total.Cells.Item(5, 1).Value2 = 3

$total.Select()
$total.Range("A1").Select()
